# Fleet resilience: hard-code number of runs.
# Each line in FleetResParam is its own experiment; every run now reuses the
# same seed (column E = 42) and the simulated run length (column H) is a
# plain literal (4380) instead of the =24*365*50 formula. The per-experiment
# run index (column Q) collapses to 1 for every row since the seed is reset
# per experiment rather than incremented across runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = 50    # A: NUM_AIRCRAFT
    $ws.Cells.Item($r, 2).Value = 50    # B: NUM_STUDENT
    $ws.Cells.Item($r, 3).Value = 40    # C: NUM_INSTRUCTOR
    $ws.Cells.Item($r, 5).Value = 42    # E: rl (same random seed for every experiment)
    $ws.Cells.Item($r, 8).Value = 4380  # H: time_line, now a literal instead of =24*365*50
    $ws.Cells.Item($r, 17).Value = 1    # Q: SLEPspots / run index, reset each experiment
}

# Match the author's resulting selection/active-cell state on the sheet.
$ws.Range("H3:H17").Select()
